# Insert a new row at position 523 (pushes existing rows 523-588 down to 524-589)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(523).Insert()

# Populate the newly inserted row 523 with the new weekly record
$ws.Cells.Item(523, 1).Value = 7
$ws.Cells.Item(523, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(523, 3).Value = "Ñuble"
$ws.Cells.Item(523, 4).Value = 45212
$ws.Cells.Item(523, 5).Value = 16
$ws.Cells.Item(523, 6).Value = "Fruta"
$ws.Cells.Item(523, 7).Value = 100104
$ws.Cells.Item(523, 8).Value = "Frutos de pepita"
$ws.Cells.Item(523, 9).Value = 100104005
$ws.Cells.Item(523, 10).Value = "Pera"
$ws.Cells.Item(523, 11).Value = "Packham's Triumph"
$ws.Cells.Item(523, 12).Value = "Primera"
$ws.Cells.Item(523, 13).Value = 180
$ws.Cells.Item(523, 14).Value = 15000
$ws.Cells.Item(523, 15).Value = 16000
$ws.Cells.Item(523, 16).Value = 15444
$ws.Cells.Item(523, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(523, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(523, 19).Value = 858
$ws.Cells.Item(523, 20).Value = 18
